# Updates the monthly "charges" (expenses) table on the active worksheet:
#  - refreshes figures for the existing rows (2-34)
#  - appends new monthly rows through the end of 2026 (rows 35-54)
#  - matches formatting of the last pre-existing data row
#  - restores the active cell selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: row, année, mois, salaires, électricité, loyer, matériels, autres
$data = @(
    @(2, 2022, "août", 0, 5000, 60000, 0, 30000),
    @(3, 2022, "septembre", 0, 5000, 60000, 0, 30000),
    @(4, 2022, "octobre", 0, 5000, 60000, 0, 30000),
    @(5, 2022, "novembre", 0, 5000, 60000, 0, 30000),
    @(6, 2022, "décembre", 0, 5000, 60000, 0, 30000),
    @(7, 2023, "janvier", 0, 5000, 60000, 0, 30000),
    @(8, 2023, "février", 0, 5000, 60000, 0, 30000),
    @(9, 2023, "mars", 0, 5000, 60000, 0, 30000),
    @(10, 2023, "avril", 0, 5000, 60000, 0, 30000),
    @(11, 2023, "mai", 0, 5000, 60000, 0, 30000),
    @(12, 2023, "juin", 0, 5000, 60000, 0, 30000),
    @(13, 2023, "juillet", 0, 5000, 60000, 0, 30000),
    @(14, 2023, "août", 0, 5000, 60000, 0, 30000),
    @(15, 2023, "septembre", 0, 5000, 60000, 0, 30000),
    @(16, 2023, "octobre", 0, 5000, 60000, 0, 30000),
    @(17, 2023, "novembre", 0, 5000, 60000, 0, 30000),
    @(18, 2023, "décembre", 0, 5000, 60000, 0, 30000),
    @(19, 2024, "janvier", 0, 5000, 60000, 0, 30000),
    @(20, 2024, "février", 0, 5000, 60000, 0, 30000),
    @(21, 2024, "mars", 1075000, 5000, 450000, 1032291, 30000),
    @(22, 2024, "avril", 1075000, 15000, 450000, 1032291, 30000),
    @(23, 2024, "mai", 1075000, 15000, 450000, 1032291, 30000),
    @(24, 2024, "juin", 1075000, 15000, 450000, 1032291, 30000),
    @(25, 2024, "juillet", 1075000, 15000, 450000, 238541, 30000),
    @(26, 2024, "août", 1225000, 15000, 450000, 246874, 30000),
    @(27, 2024, "septembre", 650000, 15000, 450000, 246874, 30000),
    @(28, 2024, "octobre", 850000, 15000, 450000, 246874, 30000),
    @(29, 2024, "novembre", 850000, 15000, 450000, 246874, 30000),
    @(30, 2024, "décembre", 850000, 15000, 450000, 246874, 30000),
    @(31, 2025, "janvier", 1155000, 15000, 450000, 246874, 30000),
    @(32, 2025, "février", 870000, 15000, 450000, 246874, 30000),
    @(33, 2025, "mars", 870000, 15000, 450000, 246874, 30000),
    @(34, 2025, "avril", 870000, 15000, 450000, 246874, 30000),
    @(35, 2025, "mai", 870000, 15000, 450000, 246874, 30000),
    @(36, 2025, "juin", 870000, 15000, 450000, 246874, 30000),
    @(37, 2025, "juillet", 870000, 15000, 450000, 246874, 30000),
    @(38, 2025, "août", 870000, 15000, 450000, 246874, 30000),
    @(39, 2025, "septembre", 870000, 15000, 450000, 246874, 30000),
    @(40, 2025, "octobre", 870000, 15000, 450000, 246874, 30000),
    @(41, 2025, "novembre", 870000, 15000, 450000, 246874, 30000),
    @(42, 2025, "décembre", 870000, 15000, 450000, 246874, 30000),
    @(43, 2026, "janvier", 870000, 15000, 450000, 246874, 30000),
    @(44, 2026, "février", 870000, 15000, 450000, 246874, 30000),
    @(45, 2026, "mars", 870000, 15000, 450000, 246874, 30000),
    @(46, 2026, "avril", 870000, 15000, 450000, 246874, 30000),
    @(47, 2026, "mai", 870000, 15000, 450000, 246874, 30000),
    @(48, 2026, "juin", 870000, 15000, 450000, 246874, 30000),
    @(49, 2026, "juillet", 870000, 15000, 450000, 246874, 30000),
    @(50, 2026, "août", 870000, 15000, 450000, 246874, 30000),
    @(51, 2026, "septembre", 870000, 15000, 450000, 246874, 30000),
    @(52, 2026, "octobre", 870000, 15000, 450000, 246874, 30000),
    @(53, 2026, "novembre", 870000, 15000, 450000, 246874, 30000),
    @(54, 2026, "décembre", 870000, 15000, 450000, 246874, 30000)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
}

# Apply the same cell formatting used on the last pre-existing row to the newly appended rows
[void]$ws.Range("A34:G34").Copy()
[void]$ws.Range("A35:G54").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Restore the active selection / scroll position
[void]$ws.Range("J50").Select()

Write-Output "Charges sheet updated: rows 2-54 refreshed."
